$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark.
#    It currently wraps an empty span right after the "decidió" run (3rd
#    bullet paragraph). The edit relocates it so it spans from the very
#    start of the document through the end of the inline picture (5th
#    paragraph) instead. Re-adding a bookmark named "_GoBack" replaces the
#    existing one (Word only ever keeps a single "_GoBack" bookmark).
# ---------------------------------------------------------------------------
$picture = $d.InlineShapes.Item(1)
$goBackRange = $d.Range(0, $picture.Range.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Fix the "travez" -> "traves" typo (inside the proofErr-wrapped run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("travez", $true, $false, $false, $false, $false, `
    $true, 1, $false, "traves", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Split the run "se comunica polimórficamente con los 3 repositorios a "
#    into two runs: "se comunica polimórficamente" and
#    " con los 3 repositorios a ". Toggling (and resetting) a character
#    formatting property on the trailing sub-range forces Word to break the
#    run at that boundary without altering the final formatting.
# ---------------------------------------------------------------------------
$splitRange = $d.Content
$splitRange.Find.Execute("se comunica polimórficamente con los 3 repositorios a ") | Out-Null

$headText = "se comunica polimórficamente"
$tailRange = $d.Range($splitRange.Start + $headText.Length, $splitRange.End)
$tailRange.Font.Bold = $true
$tailRange.Font.Bold = $false
